$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at 334, pushing all existing rows (and the historical
# data they hold) down by two - this mirrors the weekly "add this week's
# readings at the top of the block, keep history below" pattern.
$ws.Rows("334:335").Insert()

# New row 334: Primera quality entry for 2023-07-28 (date serial 45135)
$ws.Range("A334").Value = 11
$ws.Range("B334").Value = "Vega Monumental Concepción"
$ws.Range("C334").Value = "Bíobío"
$ws.Range("D334").Value = 45135
$ws.Range("E334").Value = 8
$ws.Range("F334").Value = 100114013
$ws.Range("G334").Value = "Zanahoria"
$ws.Range("H334").Value = "Sin especificar"
$ws.Range("I334").Value = "Primera"
$ws.Range("J334").Value = 500
$ws.Range("K334").Value = 5000
$ws.Range("L334").Value = 5500
$ws.Range("M334").Value = 5300
$ws.Range("N334").Value = "$/saco 20 kilos"
$ws.Range("O334").Value = "Región de Ñuble"
$ws.Range("P334").Value = 265
$ws.Range("Q334").Value = 20
$ws.Range("R334").Value = "Hortaliza"

# New row 335: Segunda quality entry for the same date (2023-07-28)
$ws.Range("A335").Value = 11
$ws.Range("B335").Value = "Vega Monumental Concepción"
$ws.Range("C335").Value = "Bíobío"
$ws.Range("D335").Value = 45135
$ws.Range("E335").Value = 8
$ws.Range("F335").Value = 100114013
$ws.Range("G335").Value = "Zanahoria"
$ws.Range("H335").Value = "Sin especificar"
$ws.Range("I335").Value = "Segunda"
$ws.Range("J335").Value = 300
$ws.Range("K335").Value = 4500
$ws.Range("L335").Value = 4500
$ws.Range("M335").Value = 4500
$ws.Range("N335").Value = "$/saco 20 kilos"
$ws.Range("O335").Value = "Región de Ñuble"
$ws.Range("P335").Value = 225
$ws.Range("Q335").Value = 20
$ws.Range("R335").Value = "Hortaliza"
